# Add a new bulleted "Google/Bing Maps API" list item right after the
# existing "Terraform" bullet at the end of the bulleted list, matching
# the same paragraph style (Prrafodelista) and numbering (numId 2).
#
# The target text is split across three runs - "Google", "/Bing" and
# " Maps API" - each sharing identical run formatting (lang=en-US), just
# like other multi-run entries already in this document (e.g.
# "Notifications" + " Hubs"). We build the new paragraph's OOXML
# directly and insert it with Range.InsertXML so the three runs are
# preserved verbatim instead of being coalesced into a single run.

$d = $word.ActiveDocument

# Locate the "Terraform" bullet - the last item of the technologies list.
$rng = $d.Content
$null = $rng.Find.Execute("Terraform", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Find.Execute narrows $rng to the matched text; walk up to its paragraph
# so we can insert right after the paragraph mark (i.e. as a new sibling
# paragraph following "Terraform").
$para = $rng.Paragraphs(1)
$insertAt = $para.Range.End

# Collapsed range positioned right after the "Terraform" paragraph mark.
$target = $d.Range($insertAt, $insertAt)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr>' +
        '<w:pStyle w:val="Prrafodelista"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
        '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Google</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>/Bing</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Maps API</w:t></w:r>' +
    '</w:p>'

$null = $target.InsertXML($newParagraphXml)
